$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F6").Value2 = 264  # was 261
$ws.Range("F7").Value2 = 13106  # was 13090
$ws.Range("F9").Value2 = 120  # was 119
$ws.Range("F10").Value2 = 273  # was 266
$ws.Range("F11").Value2 = 3387  # was 3270
$ws.Range("F13").Value2 = 6609  # was 6582
$ws.Range("F16").Value2 = 3481  # was 3470
$ws.Range("F20").Value2 = 40  # was 39
$ws.Range("F21").Value2 = 68  # was 67
$ws.Range("F22").Value2 = 124  # was 123
$ws.Range("F23").Value2 = 53  # was 51
$ws.Range("F24").Value2 = 3645  # was 3641
$ws.Range("F27").Value2 = 3174  # was 3075
$ws.Range("F29").Value2 = 1909  # was 1902
$ws.Range("F30").Value2 = 105  # was 104
$ws.Range("F31").Value2 = 230  # was 225
$ws.Range("F32").Value2 = 6782  # was 6759
$ws.Range("F34").Value2 = 1331  # was 1254
$ws.Range("F35").Value2 = 2016  # was 2009
$ws.Range("F36").Value2 = 1301  # was 1298
$ws.Range("F37").Value2 = 108  # was 106
$ws.Range("F38").Value2 = 1063  # was 1057
$ws.Range("F42").Value2 = 1152  # was 1151
$ws.Range("F43").Value2 = 1146  # was 1145
$ws.Range("F44").Value2 = 143  # was 142
$ws.Range("F45").Value2 = 1220  # was 1215
$ws.Range("F46").Value2 = 1806  # was 1800
$ws.Range("F49").Value2 = 1177  # was 1176
# Sheet: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value2 = 6  # was 5
$ws.Range("F12").Value2 = 926  # was 924
# Sheet: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value2 = 453  # was 451
$ws.Range("F3").Value2 = 622  # was 620
$ws.Range("F4").Value2 = 25  # was 23
# Sheet: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value2 = 6  # was 5
$ws.Range("F6").Value2 = 453  # was 451
$ws.Range("F7").Value2 = 622  # was 620
$ws.Range("F8").Value2 = 264  # was 261
$ws.Range("F9").Value2 = 13106  # was 13090
$ws.Range("F11").Value2 = 120  # was 119
$ws.Range("F13").Value2 = 273  # was 266
$ws.Range("F14").Value2 = 3387  # was 3270
$ws.Range("F16").Value2 = 3481  # was 3470
$ws.Range("F19").Value2 = 40  # was 39
$ws.Range("F20").Value2 = 68  # was 67
$ws.Range("F22").Value2 = 124  # was 123
$ws.Range("F23").Value2 = 53  # was 51
$ws.Range("F24").Value2 = 3645  # was 3641
$ws.Range("F27").Value2 = 3174  # was 3075
$ws.Range("F28").Value2 = 3174  # was 3076
$ws.Range("F30").Value2 = 1909  # was 1902
$ws.Range("F31").Value2 = 105  # was 104
$ws.Range("F32").Value2 = 230  # was 225
$ws.Range("F33").Value2 = 6782  # was 6759
$ws.Range("F36").Value2 = 1331  # was 1255
$ws.Range("F37").Value2 = 2016  # was 2009
$ws.Range("F39").Value2 = 1301  # was 1298
$ws.Range("F40").Value2 = 108  # was 106
$ws.Range("F41").Value2 = 1063  # was 1057
$ws.Range("F44").Value2 = 1152  # was 1151
$ws.Range("F45").Value2 = 1220  # was 1215
$ws.Range("F47").Value2 = 1806  # was 1800
